# French translation of "Facilitator guidelines - Ants Problem.docx"
#
# The document is one big table; each label/value is its own paragraph.
# We replace the *content* of specific paragraphs (identified by their
# current English text) with the French text, while leaving the
# paragraph/cell-end marks and all run formatting untouched.
#
# Direct Range.Text assignment (rather than Find.Execute replacement) is
# used so that AutoCorrect/AutoFormat "smart quote" substitution does not
# silently turn a straight apostrophe into a curly one.

$d = $word.ActiveDocument

function Set-ParagraphText($paraIndex, $oldText, $newText) {
    $p = $d.Paragraphs.Item($paraIndex)
    $r = $p.Range
    $full = $r.Text

    # Trim the trailing paragraph-mark / cell-mark control characters
    # (CR = 13, BEL/cell-mark = 7, VT = 11, ...) so we only touch the
    # visible text content of the paragraph.
    $contentLen = $full.Length
    while ($contentLen -gt 0 -and [int][char]$full[$contentLen - 1] -lt 32) {
        $contentLen -= 1
    }
    $current = $full.Substring(0, $contentLen)

    if ($current -ne $oldText) {
        throw "Paragraph $paraIndex text mismatch: expected [$oldText], found [$current]"
    }

    $sub = $d.Range($r.Start, $r.Start + $contentLen)
    $sub.Text = $newText
}

Set-ParagraphText 2  "Video Title" "Titre de la vidéo"
Set-ParagraphText 6  "Topic" "Rubrique"
Set-ParagraphText 10 "Aim(s)" "Objectif(s)"
Set-ParagraphText 14 "Length" "Durée"
Set-ParagraphText 18 "Camp Location" "Lieu du camp"
Set-ParagraphText 22 "Facilitators" "Animateurs"
Set-ParagraphText 26 "N. of students" "N. des étudiants"
Set-ParagraphText 34 "Resources" "Les ressources"
Set-ParagraphText 35 "needed" "nécessaires"
Set-ParagraphText 39 "Preparations" "Préparations"
Set-ParagraphText 45 "Video time" "Temps de la vidéo"
Set-ParagraphText 46 "What facilitator does" "Ce que fait le facilitateur"
Set-ParagraphText 47 "What learners do" "Ce que font les apprenants"
Set-ParagraphText 50 "General VMC Video Introduction" "Vidéo générale introduisant le CVM"
Set-ParagraphText 53 "00:27 – 01:08" "00:27 - 01:08"
Set-ParagraphText 54 "Video Introduction" "Video d'introduction"
Set-ParagraphText 57 "01:09 – 02:27" "01:09 - 02:27"
Set-ParagraphText 58 "Riddle" "Énigme"
Set-ParagraphText 62 "Assist the process, provoke thoughts" "Faciliter le processus, susciter des pensées"
Set-ParagraphText 66 "02:28 – 3:10" "02:28 - 3:10"
Set-ParagraphText 71 "Assist the process, provoke thoughts" "Faciliter le processus, susciter des pensées"
